$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the additional hours worked for the UML diagram week rows
$ws.Range("H13").Value = 6
$ws.Range("B14").Value = 6

# Move / update the active selection to G22 as recorded in the sheet view
$ws.Range("G22").Select()
